$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 209.5
$ws.Range("I9").Value = 127.15385
$ws.Range("K9").Value = 127.15385
$ws.Range("M9").Value = 41.84614999999999
$ws.Range("H76").Value = 13799.5
$ws.Range("I76").Value = 16566
$ws.Range("K76").Value = 16566
$ws.Range("M76").Value = -16251
$ws.Range("H79").Value = 13799.5
$ws.Range("I79").Value = 16566
$ws.Range("K79").Value = 16566
$ws.Range("M79").Value = -15474
$ws.Range("H92").Value = 3730.111
$ws.Range("I92").Value = 1566.6364
$ws.Range("K92").Value = 1566.6364
$ws.Range("M92").Value = -318.6364000000001
$ws.Range("H98").Value = 890.6
$ws.Range("I98").Value = 890.6
$ws.Range("K98").Value = 890.6
$ws.Range("M98").Value = 607.4
$ws.Range("H100").Value = 5856.294
$ws.Range("J100").Value = 8398.714
$ws.Range("L100").Value = 8398.714
$ws.Range("N100").Value = -9480.714
$ws.Range("H122").Value = 890.6
$ws.Range("I122").Value = 890.6
$ws.Range("K122").Value = 2671.8
$ws.Range("M122").Value = -221.8000000000002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 847.7143
$ws.Range("J5").Value = 1427.1428
$ws.Range("L5").Value = 1427.1428
$ws.Range("N5").Value = -1651.1428
$ws.Range("H60").Value = 73054.8
$ws.Range("I60").Value = 73054.8
$ws.Range("K60").Value = 73054.8
$ws.Range("M60").Value = -72321.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 847.7143
$ws.Range("J4").Value = 1427.1428
$ws.Range("L4").Value = 1427.1428
$ws.Range("N4").Value = -1657.1428
$ws.Range("H86").Value = 3336.6667
$ws.Range("I86").Value = 2634.5625
$ws.Range("J86").Value = 5583.4
$ws.Range("K86").Value = 2634.5625
$ws.Range("L86").Value = 5583.4
$ws.Range("M86").Value = -1511.5625
$ws.Range("N86").Value = -7829.4
$ws.Range("H89").Value = 3336.6667
$ws.Range("I89").Value = 2634.5625
$ws.Range("J89").Value = 5583.4
$ws.Range("K89").Value = 13172.8125
$ws.Range("L89").Value = 27917
$ws.Range("M89").Value = -7556.8125
$ws.Range("N89").Value = -39149
$ws.Range("H94").Value = 3063.1538
$ws.Range("I94").Value = 3552.3
$ws.Range("K94").Value = 3552.3
$ws.Range("M94").Value = -3101.3
$ws.Range("H99").Value = 2279.6
$ws.Range("I99").Value = 1849.5
$ws.Range("K99").Value = 1849.5
$ws.Range("M99").Value = -351.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41669724
$ws.Range("I31").Value = 55557930
$ws.Range("K31").Value = 55557930
$ws.Range("M31").Value = -55557635
$ws.Range("H34").Value = 41669724
$ws.Range("I34").Value = 55557930
$ws.Range("K34").Value = 55557930
$ws.Range("M34").Value = -55557728
$ws.Range("H99").Value = 13394.206
$ws.Range("I99").Value = 8949.048000000001
$ws.Range("K99").Value = 8949.048000000001
$ws.Range("M99").Value = -7451.048000000001
$ws.Range("H105").Value = 1355.6471
$ws.Range("I105").Value = 1087.0834
$ws.Range("J105").Value = 2000.2
$ws.Range("K105").Value = 1087.0834
$ws.Range("L105").Value = 2000.2
$ws.Range("M105").Value = 659.9166
$ws.Range("N105").Value = -5494.2
$ws.Range("H126").Value = 13394.206
$ws.Range("I126").Value = 8949.048000000001
$ws.Range("K126").Value = 26847.144
$ws.Range("M126").Value = -24377.144

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 12623.857
$ws.Range("I138").Value = 11030.846
$ws.Range("J138").Value = 33333
$ws.Range("K138").Value = 33092.538
$ws.Range("L138").Value = 99999
$ws.Range("M138").Value = -27952.538
$ws.Range("N138").Value = -110279

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 64997.5
$ws.Range("I40").Value = 39996
$ws.Range("K40").Value = 39996
$ws.Range("M40").Value = -39845
$ws.Range("H97").Value = 8758.267
$ws.Range("J97").Value = 12314.5
$ws.Range("L97").Value = 12314.5
$ws.Range("N97").Value = -13306.5
$ws.Range("H119").Value = 100380
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("H120").Value = 99998.336
$ws.Range("J120").Value = 99998.336
$ws.Range("L120").Value = 99998.336
$ws.Range("N120").Value = -109674.336

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 25000
$ws.Range("I4").Value = 25000
$ws.Range("K4").Value = 25000
$ws.Range("M4").Value = -24887
$ws.Range("H28").Value = 25000
$ws.Range("I28").Value = 25000
$ws.Range("K28").Value = 25000
$ws.Range("M28").Value = -24768
$ws.Range("H37").Value = 25000
$ws.Range("I37").Value = 25000
$ws.Range("K37").Value = 25000
$ws.Range("M37").Value = -24893
$ws.Range("H46").Value = 1703.9
$ws.Range("I46").Value = 1448.75
$ws.Range("J46").Value = 1874
$ws.Range("K46").Value = 1448.75
$ws.Range("L46").Value = 1874
$ws.Range("M46").Value = -1260.75
$ws.Range("N46").Value = -2250
$ws.Range("H55").Value = 1187
$ws.Range("I55").Value = 780.7143
$ws.Range("J55").Value = 1376.6
$ws.Range("K55").Value = 780.7143
$ws.Range("L55").Value = 1376.6
$ws.Range("M55").Value = -607.7143
$ws.Range("N55").Value = -1722.6
$ws.Range("H100").Value = 25030490
$ws.Range("I100").Value = 4991.8
$ws.Range("K100").Value = 4991.8
$ws.Range("M100").Value = -4450.8
$ws.Range("H132").Value = 3850.1924
$ws.Range("I132").Value = 2233.5
$ws.Range("J132").Value = 6436.9
$ws.Range("K132").Value = 6700.5
$ws.Range("L132").Value = 19310.7
$ws.Range("M132").Value = -4170.5
$ws.Range("N132").Value = -24370.7

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18348.285
$ws.Range("J45").Value = 13493.8
$ws.Range("L45").Value = 13493.8
$ws.Range("N45").Value = -14475.8
$ws.Range("H94").Value = 26713.143
$ws.Range("J94").Value = 26713.143
$ws.Range("L94").Value = 26713.143
$ws.Range("N94").Value = -28515.143
$ws.Range("H122").Value = 1816.7576
$ws.Range("I122").Value = 1432.8334
$ws.Range("J122").Value = 2277.4666
$ws.Range("K122").Value = 4298.5002
$ws.Range("L122").Value = 6832.399800000001
$ws.Range("M122").Value = -1848.5002
$ws.Range("N122").Value = -11732.3998
$ws.Range("H132").Value = 1012581.9
$ws.Range("I132").Value = 15102.375
$ws.Range("K132").Value = 45307.125
$ws.Range("M132").Value = -42777.125
$ws.Range("H136").Value = 919025.8
$ws.Range("I136").Value = 10928.4
$ws.Range("K136").Value = 32785.2
$ws.Range("M136").Value = -30235.2
